# Refresh the cryptocurrency price/volume snapshot (rows 2-51 of the only
# worksheet) to match the latest scrape. Column D ("Price") stores its values
# as text in the source file (even when they look numeric, e.g. "41.25"), so
# plain-numeric replacements are written with a leading apostrophe to stop
# Excel from silently re-typing them as numbers (which would also drop
# trailing zeros, e.g. "0.620" -> 0.62).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '42.202.68'
$ws.Range('E2').Value = '  -1.26%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '2.249.46'
$ws.Range('E3').Value = '  -1.40%  '

# Row 4: TetherUSD
$ws.Range('E4').Value = '  +0.07%  '

# Row 5: BNB
$ws.Range('D5').Value = '''247.36'
$ws.Range('E5').Value = '  -1.78%  '

# Row 6: XRP
$ws.Range('D6').Value = '''0.620'
$ws.Range('E6').Value = '  -3.78%  '

# Row 7: Solana
$ws.Range('D7').Value = '''73.86'

# Row 8: USDC
$ws.Range('E8').Value = '  -0.01%  '

# Row 9: Cardano
$ws.Range('D9').Value = '''0.612'
$ws.Range('E9').Value = '  -5.01%  '

# Row 10: Avalanche
$ws.Range('D10').Value = '''41.13'
$ws.Range('E10').Value = '  +3.64%  '

# Row 11: Dogecoin
$ws.Range('D11').Value = '''0.0934'
$ws.Range('E11').Value = '  -4.68%  '

# Row 12: Polkadot
$ws.Range('E12').Value = '  -5.36%  '

# Row 13: TRON
$ws.Range('E13').Value = '  -3.04%  '

# Row 14: WrappedliquidstakedEther2.0
$ws.Range('D14').Value = '2.584.23'
$ws.Range('E14').Value = '  -1.68%  '

# Row 15: Chainlink
$ws.Range('D15').Value = '''14.49'
$ws.Range('E15').Value = '  -3.86%  '

# Row 16: Polygon
$ws.Range('D16').Value = '''0.850'
$ws.Range('E16').Value = '  -2.50%  '

# Row 17: WrappedEther
$ws.Range('D17').Value = '2.253.72'
$ws.Range('E17').Value = '  -0.69%  '

# Row 18: WrappedBTC
$ws.Range('D18').Value = '42.106.17'
$ws.Range('E18').Value = '  -1.31%  '

# Row 19: ShibaInu
$ws.Range('D19').Value = '0.0₃0978'
$ws.Range('E19').Value = '  -2.55%  '

# Row 20: Uniswap
$ws.Range('D20').Value = '''6.12'
$ws.Range('E20').Value = '  -1.72%  '

# Row 21: Litecoin
$ws.Range('D21').Value = '''71.85'
$ws.Range('E21').Value = '  -0.74%  '

# Row 22: ImmutableX
$ws.Range('D22').Value = '''2.27'
$ws.Range('E22').Value = '  +4.84%  '

# Row 23: BitcoinCash
$ws.Range('D23').Value = '''229.96'
$ws.Range('E23').Value = '  -2.99%  '

# Row 25: Cosmos
$ws.Range('D25').Value = '''11.14'
$ws.Range('E25').Value = '  -1.61%  '

# Row 26: WEMIXToken -> InternetComputer(DFINITY)
$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D26').Value = '''7.89'
$ws.Range('E26').Value = '  +25.73%  '

# Row 27: InternetComputer(DFINITY) -> WEMIXToken
$ws.Range('B27').Value = 'WEMIXToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D27').Value = '''3.54'
$ws.Range('E27').Value = '  -8.20%  '

# Row 28: PancakeSwap
$ws.Range('E28').Value = '  -4.03%  '

# Row 29: Monero -> Toncoin
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = '''2.23'
$ws.Range('E29').Value = '  +2.78%  '

# Row 30: Toncoin -> Monero
$ws.Range('B30').Value = 'Monero'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D30').Value = '''169.29'
$ws.Range('E30').Value = '  +0.84%  '

# Row 31: EthereumClassic
$ws.Range('D31').Value = '''20.68'
$ws.Range('E31').Value = '  -1.59%  '

# Row 33: Kaspa
$ws.Range('E33').Value = '  -6.03%  '

# Row 34: InjectiveProtocol
$ws.Range('D34').Value = '''30.13'
$ws.Range('E34').Value = '  -4.26%  '

# Row 35: Stellar
$ws.Range('E35').Value = '  -2.63%  '

# Row 36: RenderToken
$ws.Range('D36').Value = '''4.49'
$ws.Range('E36').Value = '  -2.14%  '

# Row 37: Filecoin
$ws.Range('D37').Value = '''4.88'
$ws.Range('E37').Value = '  +2.19%  '

# Row 38: VeChain
$ws.Range('E38').Value = '  -1.72%  '

# Row 39: Celestia
$ws.Range('D39').Value = '''13.47'
$ws.Range('E39').Value = '  -2.07%  '

# Row 40: LidoDAOToken
$ws.Range('E40').Value = '  -5.24%  '

# Row 41: THORChain
$ws.Range('E41').Value = '  -2.17%  '

# Row 42: MultiversX
$ws.Range('D42').Value = '''62.19'
$ws.Range('E42').Value = '  +1.58%  '

# Row 43: Algorand
$ws.Range('E43').Value = '  -3.47%  '

# Row 44: Aave
$ws.Range('D44').Value = '''108.16'
$ws.Range('E44').Value = '  +2.76%  '

# Row 45: FraxShare
$ws.Range('E45').Value = '  -4.37%  '

# Row 46: Cronos
$ws.Range('E46').Value = '  -0.46%  '

# Row 47: BinanceUSD
$ws.Range('E47').Value = '  -0.25%  '

# Row 48: ARBITRUM
$ws.Range('E48').Value = '  -3.94%  '

# Row 49: TrustWalletToken
$ws.Range('E49').Value = '  -1.18%  '

# Row 50: NEARProtocol
$ws.Range('D50').Value = '''2.28'
$ws.Range('E50').Value = '  +0.61%  '

# Row 51: WOONetwork
$ws.Range('E51').Value = '  +15.43%  '
